$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts A->B, B->C, C->D, D->E, E->F
$ws.Columns.Item(1).Insert()

# Header for the new ID column, matching the style used by the other headers
# (bold font, thin border all around, centered/top-aligned)
$ws.Range("A1").Value = "ID"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4160
$ws.Range("A1").Borders.LineStyle = 1

# Row labels for the new ID column
$ids = @(
    "Hb 2",
    "Hb 3",
    "S 24",
    "S 28",
    "Hb 107",
    "Hb 66",
    "Hb 69",
    "Hb 95",
    "Hb 99",
    "Hb 92",
    "Hb 40",
    "Hb 41",
    "S 11",
    "Hb 57",
    "S 21",
    "S 22",
    "S 3",
    "S 4",
    "S 5",
    "Hb 74",
    "Hb 79",
    "Hb 32",
    "S 15",
    "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
